$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price values look numeric (single decimal point) need to be
# pre-formatted as Text so Excel stores them as strings, matching the source data
# (which uses a mix of "." as thousands separator and decimal point).
$textCells = "D5","D6","D12","D20","D21","D24","D25","D30","D32","D35","D36","D40","D41","D42","D43","D44","D47","D48","D50"
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '58.926.92'
$ws.Range("E2").Value = '  +1.09%  '
$ws.Range("D3").Value = '2.504.05'
$ws.Range("E3").Value = '  +0.78%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '535.18'
$ws.Range("E5").Value = '  +3.14%  '
$ws.Range("D6").Value = '134.30'
$ws.Range("E6").Value = '  +1.78%  '
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("D9").Value = '2.507.52'
$ws.Range("E9").Value = '  -0.40%  '
$ws.Range("E10").Value = '  +2.01%  '
$ws.Range("D12").Value = '5.17'
$ws.Range("E12").Value = '  -0.79%  '
$ws.Range("E13").Value = '  -1.54%  '
$ws.Range("D14").Value = '2.947.71'
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("D15").Value = '58.718.29'
$ws.Range("E15").Value = '  +0.87%  '
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("E17").Value = '  +0.39%  '
$ws.Range("D18").Value = '2.506.64'
$ws.Range("E18").Value = '  -0.07%  '
$ws.Range("E19").Value = '  -0.97%  '
$ws.Range("D20").Value = '4.25'
$ws.Range("E20").Value = '  +1.42%  '
$ws.Range("D21").Value = '321.18'
$ws.Range("E21").Value = '  -0.91%  '
$ws.Range("E22").Value = '  +2.91%  '
$ws.Range("E23").Value = '  +0.36%  '
$ws.Range("D24").Value = '65.81'
$ws.Range("D25").Value = '0.409'
$ws.Range("E25").Value = '  +0.34%  '
$ws.Range("E26").Value = '  +1.75%  '
$ws.Range("E27").Value = '  -1.35%  '
$ws.Range("E28").Value = '  +1.05%  '
$ws.Range("E29").Value = '  +1.00%  '
$ws.Range("D30").Value = '171.91'
$ws.Range("E30").Value = '  +2.65%  '
$ws.Range("E31").Value = '  +1.58%  '
$ws.Range("D32").Value = '6.28'
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("E33").Value = '  -0.90%  '
$ws.Range("E34").Value = '  +0.13%  '
$ws.Range("D35").Value = '0.997'
$ws.Range("E35").Value = '  +0.34%  '
$ws.Range("D36").Value = '18.09'
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("E37").Value = '  -3.72%  '
$ws.Range("E38").Value = '  +0.10%  '
$ws.Range("E39").Value = '  +3.64%  '
$ws.Range("D40").Value = '0.831'
$ws.Range("E40").Value = '  +6.33%  '
$ws.Range("D41").Value = '36.52'
$ws.Range("E41").Value = '  -0.48%  '
$ws.Range("D42").Value = '3.47'
$ws.Range("E42").Value = '  +0.93%  '
$ws.Range("D43").Value = '274.44'
$ws.Range("E43").Value = '  -1.53%  '
$ws.Range("D44").Value = '131.17'
$ws.Range("E44").Value = '  +7.34%  '
$ws.Range("E45").Value = '  -1.84%  '
$ws.Range("E46").Value = '  -1.30%  '
$ws.Range("D47").Value = '0.0936'
$ws.Range("E47").Value = '  +1.60%  '
$ws.Range("D48").Value = '0.0510'
$ws.Range("E48").Value = '  +2.16%  '
$ws.Range("E49").Value = '  +2.12%  '
$ws.Range("D50").Value = '16.82'
$ws.Range("E50").Value = '  -1.30%  '
$ws.Range("D51").Value = '1.748.90'
$ws.Range("E51").Value = '  +0.23%  '
